# Agents only migrate once every 5 years: bump the "Chance of being
# accepted (Chance of forming a bond with a male)" column on the
# Dispersal sheet from 0.66 to 0.75 for every age-class row (7..30,
# i.e. rows 9-32).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dispersal")

$ws.Range("D9:D32").Value = 0.75

# Leave the workbook focused on the sheet the author was last working
# in (Dispersal) with their last selection.
$ws.Activate()
$ws.Range("E12").Select()
